# Update "想去人数" (F column) values on the "展览" and "全部类型" worksheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new F-column value
$updates = @{
    2  = 1168
    3  = 113
    4  = 1619
    5  = 623
    8  = 11511
    12 = 361
    14 = 798
    15 = 12381
    16 = 13081
    18 = 142
    21 = 225
    23 = 45
    24 = 116
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
